$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 467; $r++) {
    $ws.Cells.Item($r, 5).Value = "V"
}
